$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 372000
$ws.Range("C3").Value = 150000

$ws.Columns.Item(1).ColumnWidth = 17.109375
$ws.Columns.Item(2).ColumnWidth = 19.109375
